$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Swap the tab names of the first two sheets ("2017 LEAVE BALANCE"
#    <-> "2018 LEAVE BALANCE"). The underlying worksheets (and their
#    data/relationship ids) stay exactly where they are - only the
#    display names trade places.
# ------------------------------------------------------------------
$wsA = $wb.Worksheets.Item(1)   # was "2017 LEAVE BALANCE"
$wsB = $wb.Worksheets.Item(2)   # was "2018 LEAVE BALANCE"

$wsA.Name = "TEMP_SWAP_NAME_1"
$wsB.Name = "2017 LEAVE BALANCE"
$wsA.Name = "2018 LEAVE BALANCE"

# Restore the Print_Titles rows (the rename above re-points them at
# the sheet's own range, so reinstate the explicit "$1:$9" rows).
$wsA.PageSetup.PrintTitleRows = "`$1:`$9"
$wsB.PageSetup.PrintTitleRows = "`$1:`$9"

# The rename can also perturb the unrelated external defined name
# BALANCE_1 on the CONVERTION sheet (it can lose the structured-
# reference part referring to the external table) - put it back.
$nExt = $wb.Names.Item("CONVERTION!BALANCE_1")
$nExt.RefersTo = "=[1]!Table1[[#Headers],[BALANCE]]"

# ------------------------------------------------------------------
# 2) Leave-card data entry on the (now) "2018 LEAVE BALANCE" sheet
#    (the worksheet that used to be named "2017 LEAVE BALANCE", still
#    worksheet #1) - a new SL(1-0-0) leave for 6/1/2023, plus filling
#    in the EARNED value for 5/1/2023 that had been left blank.
# ------------------------------------------------------------------
$ws = $wsA

# 5/1/2023 row: EARNED = 1.25 (the "EARNED " helper column recomputes
# automatically from the table's calculated formula).
$ws.Range("C83").Value = 1.25

# 6/1/2023 row: particulars, absence (w/ pay) and the leave date.
$ws.Range("B84").Value = "SL(1-0-0)"
$ws.Range("H84").Value = 1

# K84 needs the same date-formatted style already used by the other
# REMARKS-column dates (e.g. K11) - copy formatting only, then set
# the value so the shared-formula/table machinery doesn't move the
# copied value around.
$ws.Range("K11").Copy()
$ws.Range("K84").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K84").Value = 45097

# ------------------------------------------------------------------
# 3) View-state: the workbook was left scrolled to the new rows, with
#    K85 selected on the "2018 LEAVE BALANCE" sheet, and the window
#    split/scroll reset (no extra scroll) on the "2017 LEAVE BALANCE"
#    sheet.
# ------------------------------------------------------------------
$wsA.Activate()
$winA = $excel.ActiveWindow
$winA.SplitRow = 81
$ws.Range("K85").Select()

$wsB.Activate()
$winB = $excel.ActiveWindow
$winB.SplitRow = $winB.SplitRow

$wsA.Activate()
